{"js": "// Apply the day-change and the multiplication-problem updates.\n// Each entry is [oldText, newText]; every oldText is unique within the\n// document, so a single-match search+replace is safe for all of them.\nconst replacements = [\n  [\"2026-01-08 Thursday\", \"2026-01-09 Friday\"],\n  [\"27\u00d731=\", \"91\u00d792=\"],\n  [\"71\u00d751=\", \"84\u00d716=\"],\n  [\"28\u00d797=\", \"13\u00d779=\"],\n  [\"31\u00d715=\", \"34\u00d790=\"],\n  [\"95\u00d754=\", \"59\u00d731=\"],\n  [\"81\u00d792=\", \"32\u00d768=\"],\n  [\"68\u00d738=\", \"23\u00d737=\"],\n  [\"82\u00d720=\", \"42\u00d785=\"],\n  [\"23\u00d785=\", \"46\u00d740=\"],\n  [\"28\u00d730=\", \"81\u00d741=\"],\n  [\"67\u00d737=\", \"43\u00d773=\"],\n  [\"55\u00d712=\", \"86\u00d767=\"],\n  [\"69\u00d762=\", \"90\u00d711=\"],\n  [\"29\u00d754=\", \"36\u00d735=\"],\n  [\"42\u00d728=\", \"80\u00d765=\"],\n  [\"18\u00d719=\", \"97\u00d744=\"],\n  [\"69\u00d721=\", \"37\u00d756=\"],\n  [\"39\u00d745=\", \"42\u00d786=\"],\n  [\"79\u00d795=\", \"88\u00d766=\"],\n  [\"51\u00d784=\", \"45\u00d764=\"],\n  [\"96\u00d729=\", \"70\u00d758=\"],\n  [\"11\u00d756=\", \"84\u00d712=\"],\n  [\"26\u00d775=\", \"65\u00d795=\"],\n  [\"29\u00d719=\", \"97\u00d799=\"],\n  [\"62\u00d742=\", \"86\u00d796=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the day-change and the multiplication-problem updates.\n# Each pair is (oldText, newText); every oldText is unique within the\n# document, so Find/Replace All is safe for each one individually.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-08 Thursday\", \"2026-01-09 Friday\"),\n    @(\"27\u00d731=\", \"91\u00d792=\"),\n    @(\"71\u00d751=\", \"84\u00d716=\"),\n    @(\"28\u00d797=\", \"13\u00d779=\"),\n    @(\"31\u00d715=\", \"34\u00d790=\"),\n    @(\"95\u00d754=\", \"59\u00d731=\"),\n    @(\"81\u00d792=\", \"32\u00d768=\"),\n    @(\"68\u00d738=\", \"23\u00d737=\"),\n    @(\"82\u00d720=\", \"42\u00d785=\"),\n    @(\"23\u00d785=\", \"46\u00d740=\"),\n    @(\"28\u00d730=\", \"81\u00d741=\"),\n    @(\"67\u00d737=\", \"43\u00d773=\"),\n    @(\"55\u00d712=\", \"86\u00d767=\"),\n    @(\"69\u00d762=\", \"90\u00d711=\"),\n    @(\"29\u00d754=\", \"36\u00d735=\"),\n    @(\"42\u00d728=\", \"80\u00d765=\"),\n    @(\"18\u00d719=\", \"97\u00d744=\"),\n    @(\"69\u00d721=\", \"37\u00d756=\"),\n    @(\"39\u00d745=\", \"42\u00d786=\"),\n    @(\"79\u00d795=\", \"88\u00d766=\"),\n    @(\"51\u00d784=\", \"45\u00d764=\"),\n    @(\"96\u00d729=\", \"70\u00d758=\"),\n    @(\"11\u00d756=\", \"84\u00d712=\"),\n    @(\"26\u00d775=\", \"65\u00d795=\"),\n    @(\"29\u00d719=\", \"97\u00d799=\"),\n    @(\"62\u00d742=\", \"86\u00d796=\")\n)\n\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n}\n"}
